# Auto-generated script to update cryptos list (crypto prices/volumes + two row reorders)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'27.751.87"
$ws.Range("E2").Value = "  +1.33%  "

# Row 3
$ws.Range("D3").Value = "'1.878.58"
$ws.Range("E3").Value = "  +1.59%  "

# Row 4
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").Value = "'331.99"
$ws.Range("E5").Value = "  +3.56%  "

# Row 6
$ws.Range("E6").Value = "  -0.05%  "

# Row 7
$ws.Range("D7").Value = "'0.4730"
$ws.Range("E7").Value = "  +6.14%  "

# Row 8
$ws.Range("D8").Value = "'0.3961"
$ws.Range("E8").Value = "  +3.39%  "

# Row 9
$ws.Range("D9").Value = "'47.85"
$ws.Range("E9").Value = "  -0.99%  "

# Row 10
$ws.Range("D10").Value = "'0.08023"
$ws.Range("E10").Value = "  +2.82%  "

# Row 11
$ws.Range("E11").Value = "  +1.02%  "

# Row 12
$ws.Range("D12").Value = "'21.85"
$ws.Range("E12").Value = "  +1.78%  "

# Row 13
$ws.Range("D13").Value = "'1.887.59"
$ws.Range("E13").Value = "  +1.64%  "

# Row 14
$ws.Range("D14").Value = "'5.965"
$ws.Range("E14").Value = "  +2.35%  "

# Row 15
$ws.Range("D15").Value = "'7.161"
$ws.Range("E15").Value = "  +1.10%  "

# Row 16
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  -0.26%  "

# Row 17
$ws.Range("D17").Value = "'87.19"
$ws.Range("E17").Value = "  +2.08%  "

# Row 18
$ws.Range("E18").Value = "  +2.50%  "

# Row 19
$ws.Range("D19").Value = "'0.06625"
$ws.Range("E19").Value = "  +1.99%  "

# Row 20
$ws.Range("E20").Value = "  +2.38%  "

# Row 21
$ws.Range("E21").Value = "  -0.12%  "

# Row 22
$ws.Range("D22").Value = "'27.770.84"
$ws.Range("E22").Value = "  +1.41%  "

# Row 23
$ws.Range("D23").Value = "'5.511"
$ws.Range("E23").Value = "  +0.92%  "

# Row 24
$ws.Range("E24").Value = "  +2.88%  "

# Row 25
$ws.Range("E25").Value = "  +1.12%  "

# Row 26
$ws.Range("D26").Value = "'2.108.97"
$ws.Range("E26").Value = "  +1.83%  "

# Row 27
$ws.Range("D27").Value = "'156.42"
$ws.Range("E27").Value = "  +3.34%  "

# Row 28
$ws.Range("D28").Value = "'20.25"
$ws.Range("E28").Value = "  +4.97%  "

# Row 29
$ws.Range("D29").Value = "'2.099"
$ws.Range("E29").Value = "  +3.74%  "

# Row 30
$ws.Range("D30").Value = "'5.593"
$ws.Range("E30").Value = "  +2.63%  "

# Row 31
$ws.Range("D31").Value = "'122.60"
$ws.Range("E31").Value = "  +2.71%  "

# Row 32
$ws.Range("D32").Value = "'0.9688"
$ws.Range("E32").Value = "  +5.00%  "

# Row 33
$ws.Range("D33").Value = "'0.09561"
$ws.Range("E33").Value = "  +2.68%  "

# Row 34
$ws.Range("D34").Value = "'1.456"
$ws.Range("E34").Value = "  -1.84%  "

# Row 35
$ws.Range("D35").Value = "'3.626"
$ws.Range("E35").Value = "  +0.71%  "

# Row 36
$ws.Range("D36").Value = "'5.304"
$ws.Range("E36").Value = "  +2.14%  "

# Row 37
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.06116"
$ws.Range("E37").Value = "  +2.68%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02266"
$ws.Range("E38").Value = "  +2.50%  "

# Row 39
$ws.Range("D39").Value = "'1.232"
$ws.Range("E39").Value = "  +1.98%  "

# Row 40
$ws.Range("D40").Value = "'8.205"
$ws.Range("E40").Value = "  -0.79%  "

# Row 41
$ws.Range("E41").Value = "  -0.11%  "

# Row 42
$ws.Range("D42").Value = "'0.5989"
$ws.Range("E42").Value = "  +1.96%  "

# Row 43
$ws.Range("D43").Value = "'0.1912"
$ws.Range("E43").Value = "  +3.65%  "

# Row 44
$ws.Range("D44").Value = "'10.25"
$ws.Range("E44").Value = "  +0.78%  "

# Row 45
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.5717"
$ws.Range("E45").Value = "  +2.01%  "

# Row 46
$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.248"
$ws.Range("E46").Value = "  -0.36%  "

# Row 47
$ws.Range("D47").Value = "'12.29"
$ws.Range("E47").Value = "  +0.88%  "

# Row 48
$ws.Range("D48").Value = "'3.404"
$ws.Range("E48").Value = "  +1.60%  "

# Row 49
$ws.Range("D49").Value = "'1.934"
$ws.Range("E49").Value = "  +1.47%  "

# Row 50
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.00000000316"
$ws.Range("E50").Value = "  +10.73%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.06822"
$ws.Range("E51").Value = "  -0.16%  "
